# Insert a new weekly price record for "Femacal de La Calera - Espinaca"
# at row 117 (pushing the existing rows 117-180 down to 118-181), then
# populate the new row with the latest observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 117; Excel copies formatting
# (e.g. the date style on column D) from the row above automatically.
$ws.Rows.Item(117).Insert()

# Fill in the values for the newly inserted row 117.
$ws.Range("A117").Value = 3
$ws.Range("B117").Value = "Femacal de La Calera"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44452
$ws.Range("E117").Value = 5
$ws.Range("F117").Value = 100112012
$ws.Range("G117").Value = "Espinaca"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 280
$ws.Range("K117").Value = 2800
$ws.Range("L117").Value = 3000
$ws.Range("M117").Value = 2914
$ws.Range("N117").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O117").Value = "Provincia de Quillota"
$ws.Range("P117").Value = 971
$ws.Range("Q117").Value = 3
$ws.Range("R117").Value = "Hortaliza"
